# Remove one insignificant row from the student dataset.
#
# Row 49 (Roll No "214a1113" - AMBAT SHRINIVAS RAMESH MEENAKSHI, an all-zero
# placeholder record) is deleted entirely. Excel's native row-delete shifts
# every row below it up by one (old row 50 becomes the new row 49), re-points
# the shared-string-backed cell values accordingly, and drops the two shared
# strings that are no longer referenced anywhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire 49th row - everything below shifts up automatically.
$ws.Rows.Item(49).Delete()

# Mirror the resulting selection/active-cell state from the source workbook
# (Excel leaves the selection sitting on the row that used to be below the
# deleted one once the delete completes).
$ws.Activate()
$ws.Range("A49").Select()
